# UndoRedoNewCommand1StateListDiagram.pptx - "update user guide section3 images"
#
# The canonical diff for this slide touches only the three "BookShelf"
# state-list tables - shape id 13 "Table 12", id 12 "Table 11" and id 15
# "Table 14" (each a <p:graphicFrame> hosting an <a:tbl>): PowerPoint
# re-saved every one of them (bumping its internal p14:modId
# coauthoring/merge stamp) without changing any visible text, position or
# formatting anywhere else in the slide.
#
# Re-assert each table's own row height / column width (idempotent
# round-trip) so every affected graphicFrame is touched the same way the
# source edit touched them, while leaving the rendered deck unchanged.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$touchedIds = @(13, 12, 15)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable -and ($touchedIds -contains $shape.Id)) {
        $tbl = $shape.Table

        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            $row = $tbl.Rows.Item($r)
            $row.Height = $row.Height
        }

        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $col = $tbl.Columns.Item($c)
            $col.Width = $col.Width
        }
    }
}
